# Applies the "Diseñadas todas las fundaciones" edit:
#  - bumps the wall length formula on "MUROS EJE X" row 7 (E7)
#  - re-designs the footing dimensions (O/P/Q columns) for several rows of
#    the "DISEÑO" sheet, including a manually tweaked AB10 formula
#  - restores the on-screen selections / frozen-pane scroll position that
#    the author ended up with after finishing the edits

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "MUROS EJE X": wall F4X length bumped from 6.15 to 6.16 m
# ---------------------------------------------------------------------
$wsX = $wb.Worksheets.Item("MUROS EJE X")
$wsX.Range("E7").Formula = "=6.16+2*0.3"

# ---------------------------------------------------------------------
# 2. "DISEÑO": redesigned footings
# ---------------------------------------------------------------------
$wsD = $wb.Worksheets.Item("DISEÑO")

# Row 10 - wider margin (0.3 -> 0.55) and narrower footing plan (1.2 -> 1.1)
$wsD.Range("O10").Formula = "=0.55+E10+0.55"
$wsD.Range("P10").Value = 1.1
$wsD.Range("Q10").Value = 1.1
# one-off manual correction applied only to this footing's bearing check
$wsD.Range("AB10").Formula = "=MAX(IF(AA10<`$S10,(Y10/(`$P10*`$O10))-(6*Z10/(`$P10*`$O10^2)),IF(AA10=`$S10,(2*Y10)/(`$P10*`$O10),(2*Y10)/(`$P10*(3*(`$O10/2-AA10))))),IF(AA10<`$S10,(Y10/(`$P10*`$O10))+(6*Z10/(`$P10*`$O10^2)),IF(AA10=`$S10,(2*Y10)/(`$P10*`$O10),(2*Y10)/(`$P10*(3*(`$O10/2-AA10))))))/10-0.016"

# Row 30 - margin 0.5 -> 1, plan 1.6x1.4 -> 1.5x2
$wsD.Range("O30").Formula = "=1+E30+1"
$wsD.Range("P30").Value = 1.5
$wsD.Range("Q30").Value = 2

# Row 33 - plan 1x1 -> 1.2x1
$wsD.Range("P33").Value = 1.2

# Row 34 - margin 0.3 -> 0.5, plan 1x1 -> 1.5x1.3
$wsD.Range("O34").Formula = "=0.5+E34+0.5"
$wsD.Range("P34").Value = 1.5
$wsD.Range("Q34").Value = 1.3

# Row 35 - plan 1x1 -> 1.1x1
$wsD.Range("P35").Value = 1.1

# Row 36 - margin 0.3 -> 0.5, plan 1x1 -> 1.1x1
$wsD.Range("O36").Formula = "=0.5+E36+0.5"
$wsD.Range("P36").Value = 1.1

# Row 37 - plan 1x1 -> 1.5x1.3
$wsD.Range("P37").Value = 1.5
$wsD.Range("Q37").Value = 1.3

# Row 38 - plan 1x1 -> 1.1x1
$wsD.Range("P38").Value = 1.1

# Row 42 - margin 0.3 -> 0.5, plan 1x1 -> 1.2x1.1
$wsD.Range("O42").Formula = "=0.5+E42+0.5"
$wsD.Range("P42").Value = 1.2
$wsD.Range("Q42").Value = 1.1

# Row 45 - give the closing row of the table its own 2-decimal display
$wsD.Range("O45").NumberFormat = "0.00"

# ---------------------------------------------------------------------
# 3. Restore the cursor / scroll positions left behind by the author
# ---------------------------------------------------------------------
$wsX.Range("F7").Select()

$wsY = $wb.Worksheets.Item("MUROS EJE Y")
$wsY.Range("E9").Select()

$wsD.Activate()
$wsD.Range("O45").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 56
